# Updated cryptos list refresh (prices + 1h volume deltas), with two rank swaps
# (WrappedEther/Litecoin at rows 13-14, LidoDAOToken/Monero at rows 24-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new literal text. A leading "'" forces
# text storage for price strings Excel would otherwise reinterpret as
# numbers (losing trailing zeros, e.g. "5.260" -> 5.26) or that are not
# valid numbers at all (thousand-grouped prices like "27.135.70").
$updates = [ordered]@{
    "D2" = "'27.135.70"
    "E2" = "  -0.44%  "
    "D3" = "'1.890.07"
    "E3" = "  -1.02%  "
    "D4" = "1.001"
    "E4" = "  +0.11%  "
    "D5" = "306.78"
    "E5" = "  -0.50%  "
    "E6" = "  +0.15%  "
    "D7" = "0.5211"
    "E7" = "  -0.70%  "
    "D8" = "0.3753"
    "E8" = "  -0.98%  "
    "D9" = "0.07253"
    "E9" = "  -0.48%  "
    "D10" = "21.04"
    "E10" = "  -1.39%  "
    "D11" = "0.8971"
    "E11" = "  -0.43%  "
    "D12" = "0.08163"
    "E12" = "  +6.27%  "
    "B13" = "WrappedEther"
    "C13" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D13" = "'1.906.31"
    "E13" = "  -0.09%  "
    "B14" = "Litecoin"
    "C14" = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
    "D14" = "96.31"
    "E14" = "  +1.24%  "
    "D15" = "'5.260"
    "E15" = "  -0.34%  "
    "D16" = "1.003"
    "E16" = "  +0.26%  "
    "D17" = "'0.000008582"
    "E17" = "  -1.36%  "
    "D18" = "'14.50"
    "E18" = "  -0.28%  "
    "E19" = "  +0.18%  "
    "D20" = "'27.197.83"
    "E20" = "  -0.42%  "
    "D21" = "5.074"
    "E21" = "  -0.32%  "
    "D22" = "10.68"
    "E22" = "  +0.28%  "
    "D23" = "'6.390"
    "E23" = "  -0.90%  "
    "B24" = "LidoDAOToken"
    "C24" = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
    "D24" = "2.298"
    "E24" = "  -2.17%  "
    "B25" = "Monero"
    "C25" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D25" = "147.39"
    "E25" = "  +1.00%  "
    "D26" = "18.15"
    "E26" = "  -0.18%  "
    "D27" = "1.733"
    "E27" = "  -0.18%  "
    "D28" = "114.75"
    "E28" = "  -0.21%  "
    "D29" = "4.893"
    "E29" = "  -1.41%  "
    "D30" = "4.777"
    "E30" = "  -0.87%  "
    "D31" = "0.09218"
    "E31" = "  -0.24%  "
    "D32" = "0.05039"
    "E32" = "  -0.79%  "
    "D33" = "0.7897"
    "E33" = "  -0.52%  "
    "D34" = "'1.210"
    "E34" = "  -2.93%  "
    "D35" = "3.431"
    "E35" = "  +3.71%  "
    "D36" = "2.967"
    "E36" = "  -1.19%  "
    "E37" = "  -1.59%  "
    "D38" = "0.5647"
    "E38" = "  -1.06%  "
    "E39" = "  -0.73%  "
    "D40" = "1.074"
    "E40" = "  +0.01%  "
    "D41" = "8.965"
    "E41" = "  -0.63%  "
    "D42" = "6.532"
    "E42" = "  -2.09%  "
    "D43" = "115.69"
    "E43" = "  -2.91%  "
    "D44" = "0.1514"
    "E44" = "  -0.62%  "
    "D45" = "0.4845"
    "E45" = "  -0.51%  "
    "D46" = "1.001"
    "E46" = "  +0.08%  "
    "E47" = "  -1.41%  "
    "D48" = "1.614"
    "E48" = "  +0.10%  "
    "D49" = "38.02"
    "E49" = "  +1.10%  "
    "D50" = "63.14"
    "E50" = "  -1.67%  "
    "D51" = "0.05932"
    "E51" = "  -0.15%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

